$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 2286
$ws.Range("C3").Value = 10749
$ws.Range("D4").Value = 3686
$ws.Range("E5").Value = 1827
$ws.Range("F6").Value = 4063
$ws.Range("G7").Value = 14242
$ws.Range("H8").Value = 6333
$ws.Range("I9").Value = 5639
$ws.Range("J10").Value = 4325
$ws.Range("K11").Value = 1981
$ws.Range("L12").Value = 6552
$ws.Range("M13").Value = 6666
$ws.Range("N14").Value = 33567
$ws.Range("O15").Value = 4862
$ws.Range("P16").Value = 5275
$ws.Range("Q17").Value = 6831
$ws.Range("R18").Value = 5466
$ws.Range("S19").Value = 6026
$ws.Range("T20").Value = 4234
$ws.Range("U21").Value = 38641
$ws.Range("V22").Value = 8844
$ws.Range("W23").Value = 4485
$ws.Range("X24").Value = 8605
$ws.Range("Y25").Value = 5851
$ws.Range("Z26").Value = 5451
$ws.Range("AA27").Value = 5397
$ws.Range("AB28").Value = 2770
$ws.Range("AC29").Value = 3219
$ws.Range("AD30").Value = 6989
$ws.Range("AE31").Value = 3820
$ws.Range("AF32").Value = 15701
$ws.Range("AG33").Value = 8996
$ws.Range("AH34").Value = 3409
$ws.Range("AI35").Value = 4875
$ws.Range("AJ36").Value = 3326
$ws.Range("AK37").Value = 5170
$ws.Range("AL38").Value = 16401
$ws.Range("AM39").Value = 2320
$ws.Range("AN40").Value = 7957
$ws.Range("AO41").Value = 5915
$ws.Range("AP42").Value = 3926
$ws.Range("AQ43").Value = 6223
$ws.Range("AR44").Value = 12289
$ws.Range("AS45").Value = 18007
$ws.Range("AT46").Value = 2545
$ws.Range("AU47").Value = 16343
$ws.Range("AV48").Value = 5019
$ws.Range("AW49").Value = 27492
$ws.Range("AX50").Value = 2478
$ws.Range("AY51").Value = 6203
$ws.Range("AZ52").Value = 3798
$ws.Range("BA53").Value = 3874
$ws.Range("BB54").Value = 5214
$ws.Range("BC55").Value = 36404
$ws.Range("BD56").Value = 8829
$ws.Range("BE57").Value = 7789
$ws.Range("BF58").Value = 9616
$ws.Range("BG59").Value = 5474
$ws.Range("BH60").Value = 4144
$ws.Range("BI61").Value = 10136
$ws.Range("BJ62").Value = 5275
$ws.Range("BK63").Value = 4731
